$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.328.69"
$ws.Range("E2").Value = "  +6.26%  "
$ws.Range("D3").Value = "3.789.29"
$ws.Range("E3").Value = "  +22.17%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "620.69"
$ws.Range("E5").Value = "  +8.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.07"
$ws.Range("E6").Value = "  +2.27%  "
$ws.Range("D7").Value = "3.786.66"
$ws.Range("E7").Value = "  +22.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +6.09%  "
$ws.Range("E10").Value = "  +8.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.61"
$ws.Range("E11").Value = "  +3.48%  "
$ws.Range("E12").Value = "  +7.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.72"
$ws.Range("E13").Value = "  +12.61%  "
$ws.Range("E14").Value = "  +7.00%  "
$ws.Range("D15").Value = "4.411.01"
$ws.Range("E15").Value = "  +21.86%  "
$ws.Range("D16").Value = "3.776.56"
$ws.Range("E16").Value = "  +21.64%  "
$ws.Range("D17").Value = "71.409.01"
$ws.Range("E17").Value = "  +6.45%  "
$ws.Range("E18").Value = "  +1.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.58"
$ws.Range("E19").Value = "  +7.90%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.98"
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "521.03"
$ws.Range("E21").Value = "  +6.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.44"
$ws.Range("E22").Value = "  +22.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.753"
$ws.Range("E23").Value = "  +9.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.55"
$ws.Range("E24").Value = "  +12.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "88.96"
$ws.Range("E25").Value = "  +6.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.62"
$ws.Range("E26").Value = "  +8.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.18"
$ws.Range("E27").Value = "  +9.87%  "
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("E29").Value = "  +11.02%  "
$ws.Range("E30").Value = "  +3.49%  "
$ws.Range("E31").Value = "  +12.24%  "
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0000113"
$ws.Range("E32").Value = "  +20.47%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "32.35"
$ws.Range("E33").Value = "  +15.18%  "
$ws.Range("E34").Value = "  +5.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("E36").Value = "  +11.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.17"
$ws.Range("E37").Value = "  +10.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.345"
$ws.Range("E38").Value = "  +10.30%  "
$ws.Range("E39").Value = "  +10.24%  "
$ws.Range("E40").Value = "  +9.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.86"
$ws.Range("E41").Value = "  +5.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "443.02"
$ws.Range("E42").Value = "  +20.04%  "
$ws.Range("D43").Value = "3.179.83"
$ws.Range("E43").Value = "  +13.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "45.03"
$ws.Range("E44").Value = "  -4.94%  "
$ws.Range("E45").Value = "  +7.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.84"
$ws.Range("E46").Value = "  +5.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0367"
$ws.Range("E47").Value = "  +6.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "28.17"
$ws.Range("E48").Value = "  +10.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.40"
$ws.Range("E49").Value = "  +3.56%  "
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.48"
$ws.Range("E51").Value = "  +8.60%  "
